# Updates license on the footer credit line from "CC BY" to "CC BY SA".
#
# Slide 1: the footer textbox is cut from its original spot (right after
# the "readr.png" picture) and pasted back, which PowerPoint appends at
# the end of the shape tree -- matching the target layout where the
# (re-created) footer shape is now the very last shape on the slide.
# Slide 3: the footer textbox is edited in place (no reordering needed).

$p = $ppt.ActivePresentation

function Update-RStudioFooter($shape) {
    $tr = $shape.TextFrame.TextRange
    $fullText = $tr.Text

    $ccIdx = $fullText.IndexOf("CC BY ")
    if ($ccIdx -ge 0) {
        # 1-based character index for TextRange.Characters
        $start = $ccIdx + 1
        $ccRun = $tr.Characters($start, 6)
        $ccRun.Text = "CC BY SA"

        $afterStart = $start + 8
        $nextRun = $tr.Characters($afterStart, 11)
        $nextRun.Text = "  RStudio •  "
    }

    # Re-assert the autofit height PowerPoint recalculated while we were
    # editing the run text, so the shape keeps its original extent.
    $shape.Height = 19.59386

    $newName = $shape.Name
    $newName = $newName.Replace("CC BY RStudio", "CC BY SA  RStudio")
    $shape.Name = $newName
}

# ---- Slide 1 -----------------------------------------------------------
$s1 = $p.Slides.Item(1)
$footer1 = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $cand = $s1.Shapes.Item($i)
    if ($cand.Name.StartsWith("RStudio")) {
        $footer1 = $cand
        break
    }
}
$footer1.Cut()
$pastedRange = $s1.Shapes.Paste()
$footer1New = $s1.Shapes.Item($s1.Shapes.Count)
Update-RStudioFooter $footer1New

# ---- Slide 3 ------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$footer3 = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $cand = $s3.Shapes.Item($i)
    if ($cand.Name.StartsWith("RStudio")) {
        $footer3 = $cand
        break
    }
}
Update-RStudioFooter $footer3
